# Weekly update: add two new daily price rows (new reporting date 44785)
# for "Agrícola del Norte S.A. de Arica - Locoto" before the existing
# rows 105.. (which all shift down by two rows to 107..123).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 105 - existing data (old rows
# 105-121) shifts down to 107-123.
$ws.Rows.Item(105).Insert()
$ws.Rows.Item(105).Insert()

# New row 105: Primera quality, new reporting date.
$ws.Cells.Item(105, 1).Value = 1
$ws.Cells.Item(105, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(105, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(105, 4).Value = 44785
$ws.Cells.Item(105, 5).Value = 15
$ws.Cells.Item(105, 6).Value = 100112042
$ws.Cells.Item(105, 7).Value = "Locoto"
$ws.Cells.Item(105, 8).Value = "Sin especificar"
$ws.Cells.Item(105, 9).Value = "Primera"
$ws.Cells.Item(105, 10).Value = 130
$ws.Cells.Item(105, 11).Value = 24000
$ws.Cells.Item(105, 12).Value = 25000
$ws.Cells.Item(105, 13).Value = 24500
$ws.Cells.Item(105, 14).Value = "$/caja 20 kilos"
$ws.Cells.Item(105, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(105, 16).Value = 1225
$ws.Cells.Item(105, 17).Value = 20
$ws.Cells.Item(105, 18).Value = "Hortaliza"

# New row 106: Segunda quality, same reporting date.
$ws.Cells.Item(106, 1).Value = 1
$ws.Cells.Item(106, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(106, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(106, 4).Value = 44785
$ws.Cells.Item(106, 5).Value = 15
$ws.Cells.Item(106, 6).Value = 100112042
$ws.Cells.Item(106, 7).Value = "Locoto"
$ws.Cells.Item(106, 8).Value = "Sin especificar"
$ws.Cells.Item(106, 9).Value = "Segunda"
$ws.Cells.Item(106, 10).Value = 120
$ws.Cells.Item(106, 11).Value = 19000
$ws.Cells.Item(106, 12).Value = 20000
$ws.Cells.Item(106, 13).Value = 19500
$ws.Cells.Item(106, 14).Value = "$/caja 20 kilos"
$ws.Cells.Item(106, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(106, 16).Value = 975
$ws.Cells.Item(106, 17).Value = 20
$ws.Cells.Item(106, 18).Value = "Hortaliza"
